# DateFormatTests.xlsx — add four new "Time" category test rows (43-46) to the
# "Tests" sheet that exercise TEXT() formats with an escaped-literal "days"
# token (\d\a\y\s) vs. a quoted-literal token ("days"), to make sure the
# hour token position is handled correctly (github-234).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tests")

# Serial date/time value reused by the other "Time" category rows already in
# the sheet (row 42): 1952-10-11 14:35:27 under date1904.
$timeValue = 17816.607951388887

# ---- Row 43: d \d\a\y\s h -----------------------------------------------
$ws.Range("B43").Value = "d \d\a\y\s h"
$ws.Range("B43").Font.Name = "Lucida Sans Regular"
$ws.Range("B43").NumberFormat = "@"
$ws.Range("C43").Value = $timeValue
$ws.Range("C43").NumberFormat = "dd\-mmm\-yyyy\ hh:mm:ss.000"
$ws.Range("D43").Value = "Time"

# ---- Row 44: d "days" h ---------------------------------------------------
$ws.Range("B44").Value = 'd "days" h'
$ws.Range("B44").Font.Name = "Lucida Sans Regular"
$ws.Range("B44").NumberFormat = "@"
$ws.Range("C44").Value = $timeValue
$ws.Range("C44").NumberFormat = "dd\-mmm\-yyyy\ hh:mm:ss.000"
$ws.Range("D44").Value = "Time"

# ---- Row 45: d \d\a\y\s h a/p ---------------------------------------------
$ws.Range("B45").Value = "d \d\a\y\s h a/p"
$ws.Range("B45").Font.Name = "Lucida Sans Regular"
$ws.Range("B45").NumberFormat = "@"
$ws.Range("C45").Value = $timeValue
$ws.Range("C45").NumberFormat = "dd\-mmm\-yyyy\ hh:mm:ss.000"
$ws.Range("D45").Value = "Time"

# ---- Row 46: d "days" h am/pm ----------------------------------------------
$ws.Range("B46").Value = 'd "days" h am/pm'
$ws.Range("B46").Font.Name = "Lucida Sans Regular"
$ws.Range("B46").NumberFormat = "@"
$ws.Range("C46").Value = $timeValue
$ws.Range("C46").NumberFormat = "dd\-mmm\-yyyy\ hh:mm:ss.000"
$ws.Range("D46").Value = "Time"

# Fill the TEXT(C,B) formula down from the existing shared formula (A40:A42)
# through the four new rows so each evaluates relative to its own row.
$ws.Range("A40:A46").Formula = "=TEXT(C40,B40)"

# Widen column C (format-code column) now that it holds longer strings, and
# drop the old best-fit autosize in favor of an explicit width.
$ws.Columns.Item(3).ColumnWidth = 48.75

# Keep the sheet's active selection/view in sync with the last data cell,
# matching a manual edit ending at the bottom of the new rows.
$ws.Activate()
$ws.Range("C46").Select()

Write-Host "A43:" $ws.Range("A43").Text
Write-Host "A44:" $ws.Range("A44").Text
Write-Host "A45:" $ws.Range("A45").Text
Write-Host "A46:" $ws.Range("A46").Text
